$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "In Translation" -------------------
# Overview sheet: zh-cn / de-de status columns (E,F) for both data rows
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# zh-cn sheet: Status column (C) for both data rows
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

# de-de sheet: Status column (C) for both data rows
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrower status-related columns (report generated for archive) --------
# Target stored width ~13.41 chars (was ~17.22 chars). ColumnWidth is quantized
# by the host to a 1/6-character grid, so feed it the char value whose
# snapped result lands closest to the target.
$overview.Columns.Item(5).ColumnWidth = 12.5   # -> stored width 13.3333...
$overview.Columns.Item(6).ColumnWidth = 12.5   # -> stored width 13.3333...

$zhcn.Columns.Item(3).ColumnWidth = 12.5        # -> stored width 13.3333...

$dede.Columns.Item(3).ColumnWidth = 12.5        # -> stored width 13.3333...
